# Update "想去人数" (wanted-to-go count) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览" - row => new value for column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1326
$ws1.Range("F4").Value  = 14575
$ws1.Range("F5").Value  = 17572
$ws1.Range("F16").Value = 35
$ws1.Range("F18").Value = 42
$ws1.Range("F19").Value = 1325
$ws1.Range("F24").Value = 7218
$ws1.Range("F30").Value = 5857
$ws1.Range("F32").Value = 44
$ws1.Range("F33").Value = 134
$ws1.Range("F36").Value = 5077

# Sheet "全部类型" - row => new value for column F
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 1326
$ws4.Range("F4").Value  = 14575
$ws4.Range("F5").Value  = 17572
$ws4.Range("F16").Value = 35
$ws4.Range("F18").Value = 42
$ws4.Range("F19").Value = 1325
$ws4.Range("F25").Value = 7218
$ws4.Range("F32").Value = 5857
$ws4.Range("F34").Value = 44
$ws4.Range("F35").Value = 134
$ws4.Range("F38").Value = 5077
